$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update point_id column (B) for rows 3-28 to be sequential values 2..27
for ($row = 3; $row -le 28; $row++) {
    $ws.Cells.Item($row, 2).Value = $row - 1
}

# Reflect the new view/selection state: topLeftCell moved to A3 and the
# active cell/selection moved to F24.
$ws.Activate()
$ws.Range("F24").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
